$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''56.158.15'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '''  -3.17%  '
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = '''2.367.46'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '''  -3.49%  '
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = '''  -0.04%  '
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = '''501.55'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '''  -1.80%  '
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = '''129.07'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '''  -3.49%  '
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = '''  -0.14%  '
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = '''0.545'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '''  -2.29%  '
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = '''2.371.28'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '''  -3.31%  '
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = '''0.0982'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '''  +0.14%  '
$ws.Range("E10").ClearFormats()
$ws.Range("E11").Value = '''  +0.30%  '
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = '''4.88'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '''  +5.72%  '
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = '''0.323'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '''  -0.38%  '
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = '''2.788.14'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '''  -3.49%  '
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = '''56.115.08'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '''  -3.03%  '
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = '''21.43'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '''  -2.32%  '
$ws.Range("E16").ClearFormats()
$ws.Range("E17").Value = '''  -1.69%  '
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = '''2.352.35'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '''  -4.57%  '
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = '''10.00'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '''  -3.20%  '
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = '''4.04'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '''  -2.82%  '
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = '''307.02'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '''  -2.46%  '
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = '''6.27'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '''  -2.65%  '
$ws.Range("E22").ClearFormats()
$ws.Range("E23").Value = '''  -0.12%  '
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = '''65.81'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '''  +0.53%  '
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = '''0.998'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '''  -0.12%  '
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = '''0.370'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '''  -2.69%  '
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = '''0.147'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '''  -5.78%  '
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = '''7.21'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '''  -4.77%  '
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = '''170.99'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '''  -0.88%  '
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = '''0.0₃0711'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '''  -2.96%  '
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = '''1.64'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '''  -3.47%  '
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = '''  +0.04%  '
$ws.Range("E32").ClearFormats()
$ws.Range("B33").Value = '''Aptos'
$ws.Range("B33").ClearFormats()
$ws.Range("C33").Value = '''https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("C33").ClearFormats()
$ws.Range("D33").Value = '''5.75'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '''  -6.71%  '
$ws.Range("E33").ClearFormats()
$ws.Range("E34").Value = '''  -4.47%  '
$ws.Range("E34").ClearFormats()
$ws.Range("B35").Value = '''FirstDigitalUSD'
$ws.Range("B35").ClearFormats()
$ws.Range("C35").Value = '''https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("C35").ClearFormats()
$ws.Range("D35").Value = '''0.997'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '''  -0.21%  '
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = '''17.56'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '''  -2.94%  '
$ws.Range("E36").ClearFormats()
$ws.Range("E37").Value = '''  -6.26%  '
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = '''3.74'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '''  -3.12%  '
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = '''36.07'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '''  -1.80%  '
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = '''0.791'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '''  -1.69%  '
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = '''  -5.63%  '
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = '''129.22'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '''  -5.59%  '
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = '''3.35'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '''  -1.58%  '
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = '''4.69'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '''  -4.17%  '
$ws.Range("E44").ClearFormats()
$ws.Range("E45").Value = '''  -2.32%  '
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = '''0.0901'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '''  -1.97%  '
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = '''239.35'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '''  -6.71%  '
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = '''0.0480'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '''  -2.54%  '
$ws.Range("E48").ClearFormats()
$ws.Range("E49").Value = '''  -3.96%  '
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = '''17.05'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '''  -0.70%  '
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = '''0.950'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '''  -0.66%  '
$ws.Range("E51").ClearFormats()
